$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 3 de Octubre de 2020 a las 06:59"

# India (row 5) - updated case counts
$ws.Range("B5").Value = 6473544
$ws.Range("C5").Value = 1610
$ws.Range("D5").Value = 5427706
$ws.Range("E5").Value = 944963

# Uzbekistan (row 60) - updated case counts
$ws.Range("B60").Value = 57776
$ws.Range("C60").Value = 322
$ws.Range("D60").Value = 54456
$ws.Range("E60").Value = 2845
$ws.Range("G60").Value = 3
$ws.Range("H60").Value = 475

# Tailandia (row 141) - updated case counts
$ws.Range("B141").Value = 3583
$ws.Range("C141").Value = 8
$ws.Range("D141").Value = 3386
$ws.Range("E141").Value = 138

# Belice now overtakes Yemen in ranking, so row 156 becomes Belice (with new
# numbers) and row 157 becomes Yemen (keeping its former row-156 numbers).
$ws.Range("A156").Value = "Belice"
$ws.Range("B156").Value = 2080
$ws.Range("C156").Value = 54
$ws.Range("D156").Value = 1290
$ws.Range("E156").Value = 762
$ws.Range("G156").Value = 1
$ws.Range("H156").Value = 28

$ws.Range("A157").Value = "Yemen"
$ws.Range("B157").Value = 2040
$ws.Range("C157").Value = 0
$ws.Range("D157").Value = 1307
$ws.Range("E157").Value = 144
$ws.Range("G157").Value = 0
$ws.Range("H157").Value = 589

# Butan (row 187) - updated case counts
$ws.Range("B187").Value = 283
$ws.Range("C187").Value = 1
$ws.Range("D187").Value = 229
$ws.Range("E187").Value = 54
